$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update row 10's activity description (B10).
$ws.Range("B10").Value = "Izveidot izmēģinājuma pārbaudes darbus testēšanai, prasīt skolēniem tos aizpildīt un noskenēt un nobildēt un noteikt vai "

# Widen column B to fit the longer text (accounts for the runtime's pixel
# quantization of ColumnWidth; this input value lands on a stored width of
# 105, the closest achievable value to the target 105.109375).
$ws.Columns.Item(2).ColumnWidth = 104.2

# Move the active selection to E23.
$ws.Range("E23").Select()
